$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerate the "K" column (column G) values for each saved game row.
$kValues = @{
    2  = 4
    3  = 3
    4  = 2
    5  = 3
    6  = 6
    7  = 6
    8  = 2
    9  = 1
    10 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
